$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. A leading "'" forces Excel to store the value
# as text (matching the source file's inlineStr cells) instead of auto-
# converting number-like strings (e.g. "216.80") into numeric values.
$updates = [ordered]@{
    "D2" = '''26.869.80'
    "E2" = '''  +0.30%  '
    "D3" = '''1.637.47'
    "E3" = '''  -0.18%  '
    "E4" = '''  -0.76%  '
    "D5" = '''216.80'
    "E5" = '''  -0.76%  '
    "E6" = '''  +1.83%  '
    "E7" = '''  -0.68%  '
    "E8" = '''  +1.55%  '
    "E9" = '''  +0.72%  '
    "D10" = '''19.95'
    "E10" = '''  +3.95%  '
    "D11" = '''0.0847'
    "E11" = '''  +0.10%  '
    "D12" = '''1.866.52'
    "E12" = '''  -0.17%  '
    "D13" = '''1.620.22'
    "E13" = '''  -1.19%  '
    "E14" = '''  -0.81%  '
    "D15" = '''0.530'
    "E15" = '''  +0.82%  '
    "D16" = '''66.91'
    "E16" = '''  +2.97%  '
    "D17" = '''26.854.40'
    "E17" = '''  +0.20%  '
    "E18" = '''  -0.15%  '
    "D19" = '''219.79'
    "E19" = '''  +1.66%  '
    "E20" = '''  -0.62%  '
    "D21" = '''6.80'
    "E21" = '''  +2.84%  '
    "E23" = '''  +3.96%  '
    "E24" = '''  +0.13%  '
    "D25" = '''146.62'
    "E25" = '''  -0.61%  '
    "E26" = '''  -0.65%  '
    "D27" = '''7.34'
    "E27" = '''  +3.43%  '
    "E28" = '''  +0.85%  '
    "D29" = '''15.78'
    "E29" = '''  +0.46%  '
    "D30" = '''0.0504'
    "E30" = '''  -0.50%  '
    "E31" = '''  -1.24%  '
    "E32" = '''  -1.67%  '
    "E33" = '''  +0.73%  '
    "E34" = '''  +0.87%  '
    "D35" = '''1.255.58'
    "E35" = '''  -0.69%  '
    "E36" = '''  -0.43%  '
    "E37" = '''  +2.01%  '
    "D38" = '''0.535'
    "E38" = '''  +0.89%  '
    "D39" = '''0.832'
    "E40" = '''  -0.63%  '
    "D41" = '''0.810'
    "E41" = '''  +0.49%  '
    "D42" = '''5.39'
    "E42" = '''  +0.92%  '
    "D43" = '''1.776.76'
    "E43" = '''  -0.20%  '
    "B44" = '''Aave'
    "C44" = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    "D44" = '''61.74'
    "E44" = '''  +1.80%  '
    "B45" = '''MXToken'
    "C45" = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    "D45" = '''2.10'
    "E45" = '''  -2.15%  '
    "D46" = '''91.57'
    "E46" = '''  -0.78%  '
    "E47" = '''  -0.88%  '
    "E48" = '''  +3.00%  '
    "E49" = '''  -0.40%  '
    "D50" = '''7.65'
    "E50" = '''  +1.25%  '
    "E51" = '''  -0.17%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = $updates[$addr]
    # Drop the quote-prefix style Excel applies for text-forced numeric
    # strings so the cell format matches the original (no explicit style).
    $cell.Style = "Normal"
}
